$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update price/volume(1h) columns per latest scrape.
# Numeric-looking price strings need an explicit Text format so Excel
# keeps them as literal strings (matching the source scraper output)
# instead of re-parsing them into floating point numbers.
$ws.Range("D2").Value = "56.075.90"
$ws.Range("E2").Value = "  -3.38%  "
$ws.Range("D3").Value = "2.364.38"
$ws.Range("E3").Value = "  -3.66%  "
$ws.Range("E4").Value = "  +0.18%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "500.55"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "128.73"
$ws.Range("E6").Value = "  -3.96%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -2.25%  "
$ws.Range("D9").Value = "2.366.72"
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("E11").Value = "  +0.24%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "4.86"
$ws.Range("E12").Value = "  +5.13%  "
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "2.783.70"
$ws.Range("E14").Value = "  -3.72%  "
$ws.Range("D15").Value = "56.048.63"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").Value = "2.350.36"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("E20").Value = "  -2.50%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "306.55"
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("E23").Value = "  -0.14%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "65.96"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("E27").Value = "  -6.14%  "
$ws.Range("E28").Value = "  -4.91%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "171.21"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "0.0₃0708"
$ws.Range("E30").Value = "  -3.42%  "
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("E32").Value = "  +0.19%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.997"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.08"
$ws.Range("E34").Value = "  -4.89%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.73"
$ws.Range("E35").Value = "  -7.30%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "17.59"
$ws.Range("E36").Value = "  -2.64%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.17"
$ws.Range("E37").Value = "  -6.49%  "
$ws.Range("E38").Value = "  -3.71%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "36.06"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("E40").Value = "  -2.57%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.38"
$ws.Range("E41").Value = "  -5.92%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "128.85"
$ws.Range("E42").Value = "  -5.96%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.34"
$ws.Range("E43").Value = "  -1.75%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "4.67"
$ws.Range("E44").Value = "  -4.70%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.564"
$ws.Range("E45").Value = "  -2.31%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0900"
$ws.Range("E46").Value = "  -1.93%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "239.21"
$ws.Range("E47").Value = "  -6.89%  "
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("E49").Value = "  -3.86%  "
$ws.Range("E50").Value = "  -1.27%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.949"
$ws.Range("E51").Value = "  -0.83%  "
